# "All done with all Example reconciliations."
# Update the two column headers on the Summary sheet: drop the leading
# "Original # of " from each label, then tighten columns C:D to the new
# (shorter) best-fit widths, and finally leave the selection on H20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Update header text (C7 / D7) -----------------------------------------
$ws.Range("C7").Value = "Proposed SME labels (includes duplicates)"
$ws.Range("D7").Value = "Proposed SME labels (without duplicates)"

# --- Re-fit columns C and D to their new (shorter) text --------------------
# Column C and D used to share one merged <col min="3" max="4"> width of
# 47.33203125. With the shorter text they now best-fit to two distinct,
# narrower widths (~36.16 and ~35.83 characters respectively).
$ws.Columns.Item(3).ColumnWidth = 212 / 6
$ws.Columns.Item(4).ColumnWidth = 210 / 6

# --- Leave the selection where the author left it --------------------------
$ws.Range("H20").Select() | Out-Null
